# Update existing row 2: the email in A2 picked up a typo ("comm"); leave the
# existing hyperlink (rId1 -> mailto:shahid+257@troontechnologies.com) as-is,
# matching how Excel treats hyperlink metadata as independent of cell text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "shahid+257@troontechnologies.comm"

# Row 3: valid email / bad password (typo'd with trailing #)
$ws.Range("A3").Value = "shahid+257@troontechnologies.com"
$ws.Range("B3").Value = "12345Qwe!@#"

# Row 4 is intentionally left blank (skipped)

# Row 5: repeat of row 2's pair
$ws.Range("A5").Value = "shahid+257@troontechnologies.comm"
$ws.Range("B5").Value = "12345Qwe!@"

# Row 6: valid email / valid password
$ws.Range("A6").Value = "shahid+257@troontechnologies.com"
$ws.Range("B6").Value = "12345Qwe!@"

# Wire up hyperlinks for the new cells, in the same left-to-right,
# top-to-bottom order Excel would assign them (A3, A5, A6, B3, B5, B6).
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:shahid+257@troontechnologies.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:shahid+257@troontechnologies.comm") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:shahid+257@troontechnologies.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:12345Qwe!@#") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:12345Qwe!@") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:12345Qwe!@") | Out-Null

# Hyperlinks.Add stamps a brand-new "hyperlink" cell style the first time it
# runs; bring every new cell back to the same shared Hyperlink style (s="1")
# that the pre-existing linked cells already use.
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("A5").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("A6").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"

$ws.Range("B6").Select() | Out-Null
